$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 1673.7916
$ws.Range("I33").Value = 1918.1052
$ws.Range("K33").Value = 1918.1052
$ws.Range("M33").Value = -1689.1052
$ws.Range("H43").Value = 1186.4445
$ws.Range("J43").Value = 1561.3334
$ws.Range("L43").Value = 1561.3334
$ws.Range("N43").Value = -1699.3334
$ws.Range("H51").Value = 8937.68
$ws.Range("J51").Value = 9258.869000000001
$ws.Range("L51").Value = 9258.869000000001
$ws.Range("N51").Value = -10226.869
$ws.Range("H88").Value = 28398.312
$ws.Range("I88").Value = 2212.8572
$ws.Range("J88").Value = 48764.777
$ws.Range("K88").Value = 2212.8572
$ws.Range("L88").Value = 48764.777
$ws.Range("M88").Value = -1806.8572
$ws.Range("N88").Value = -49576.777
$ws.Range("H91").Value = 28398.312
$ws.Range("I91").Value = 2212.8572
$ws.Range("J91").Value = 48764.777
$ws.Range("K91").Value = 2212.8572
$ws.Range("L91").Value = 48764.777
$ws.Range("M91").Value = -808.8571999999999
$ws.Range("N91").Value = -51572.777

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1351.5625
$ws.Range("I2").Value = 1315.5676
$ws.Range("J2").Value = 1472.6364
$ws.Range("K2").Value = 1315.5676
$ws.Range("L2").Value = 1472.6364
$ws.Range("M2").Value = -1202.5676
$ws.Range("N2").Value = -1698.6364
$ws.Range("H45").Value = 1245.0385
$ws.Range("I45").Value = 1151.4762
$ws.Range("K45").Value = 1151.4762
$ws.Range("M45").Value = -774.4762000000001
$ws.Range("H116").Value = 1351.5625
$ws.Range("I116").Value = 1315.5676
$ws.Range("J116").Value = 1472.6364
$ws.Range("K116").Value = 1315.5676
$ws.Range("L116").Value = 1472.6364
$ws.Range("M116").Value = 978.4323999999999
$ws.Range("N116").Value = -6060.6364
$ws.Range("H132").Value = 3280.7026
$ws.Range("I132").Value = 2779.7
$ws.Range("K132").Value = 8339.099999999999
$ws.Range("M132").Value = -5809.099999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1351.5625
$ws.Range("I3").Value = 1315.5676
$ws.Range("J3").Value = 1472.6364
$ws.Range("K3").Value = 1315.5676
$ws.Range("L3").Value = 1472.6364
$ws.Range("M3").Value = -1201.5676
$ws.Range("N3").Value = -1700.6364
$ws.Range("H52").Value = 83863.875
$ws.Range("J52").Value = 49898.8
$ws.Range("L52").Value = 49898.8
$ws.Range("N52").Value = -50424.8
$ws.Range("H94").Value = 1161.4062
$ws.Range("I94").Value = 738.6
$ws.Range("K94").Value = 738.6
$ws.Range("M94").Value = -287.6
$ws.Range("H121").Value = 83863.875
$ws.Range("J121").Value = 49898.8
$ws.Range("L121").Value = 49898.8
$ws.Range("N121").Value = -53392.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2917.9656
$ws.Range("J31").Value = 5012.5293
$ws.Range("L31").Value = 5012.5293
$ws.Range("N31").Value = -5602.5293
$ws.Range("H34").Value = 2917.9656
$ws.Range("J34").Value = 5012.5293
$ws.Range("L34").Value = 5012.5293
$ws.Range("N34").Value = -5416.5293
$ws.Range("H86").Value = 3334.75
$ws.Range("I86").Value = 3525
$ws.Range("J86").Value = 2764
$ws.Range("K86").Value = 3525
$ws.Range("L86").Value = 2764
$ws.Range("M86").Value = -2402
$ws.Range("N86").Value = -5010
$ws.Range("H89").Value = 3334.75
$ws.Range("I89").Value = 3525
$ws.Range("J89").Value = 2764
$ws.Range("K89").Value = 17625
$ws.Range("L89").Value = 13820
$ws.Range("M89").Value = -12009
$ws.Range("N89").Value = -25052
$ws.Range("H99").Value = 7717.0454
$ws.Range("I99").Value = 6361.3125
$ws.Range("K99").Value = 6361.3125
$ws.Range("M99").Value = -4863.3125
$ws.Range("H107").Value = 1160.8813
$ws.Range("I107").Value = 637.25
$ws.Range("J107").Value = 2696.8667
$ws.Range("K107").Value = 637.25
$ws.Range("L107").Value = 2696.8667
$ws.Range("M107").Value = 1282.75
$ws.Range("N107").Value = -6536.8667
$ws.Range("H126").Value = 7717.0454
$ws.Range("I126").Value = 6361.3125
$ws.Range("K126").Value = 19083.9375
$ws.Range("M126").Value = -16613.9375
$ws.Range("H132").Value = 5961.909
$ws.Range("I132").Value = 3389.5806
$ws.Range("J132").Value = 45833
$ws.Range("K132").Value = 10168.7418
$ws.Range("L132").Value = 137499
$ws.Range("M132").Value = -7638.7418
$ws.Range("N132").Value = -142559

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 298.5
$ws.Range("J34").Value = 297.5
$ws.Range("L34").Value = 892.5
$ws.Range("N34").Value = -1060.5
$ws.Range("H93").Value = 9512.5
$ws.Range("J93").Value = 9512.5
$ws.Range("L93").Value = 28537.5
$ws.Range("N93").Value = -32281.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H41").Value = 1800
$ws.Range("I41").Value = 1800
$ws.Range("K41").Value = 1800
$ws.Range("M41").Value = -1445
$ws.Range("H80").Value = 8944.344999999999
$ws.Range("I80").Value = 4239.6
$ws.Range("K80").Value = 4239.6
$ws.Range("M80").Value = -3241.6
$ws.Range("H83").Value = 8944.344999999999
$ws.Range("I83").Value = 4239.6
$ws.Range("K83").Value = 21198
$ws.Range("M83").Value = -16206
$ws.Range("H102").Value = 20895.254
$ws.Range("I102").Value = 24263.627
$ws.Range("K102").Value = 24263.627
$ws.Range("M102").Value = -22641.627
$ws.Range("H132").Value = 2491.3784
$ws.Range("I132").Value = 2142.0967
$ws.Range("J132").Value = 4296
$ws.Range("K132").Value = 6426.2901
$ws.Range("L132").Value = 12888
$ws.Range("M132").Value = -3896.2901
$ws.Range("N132").Value = -17948

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 47298.15
$ws.Range("J7").Value = 2252
$ws.Range("L7").Value = 2252
$ws.Range("N7").Value = -2476
$ws.Range("H22").Value = 449
$ws.Range("I22").Value = 200
$ws.Range("K22").Value = 200
$ws.Range("M22").Value = 95
$ws.Range("H27").Value = 449
$ws.Range("I27").Value = 200
$ws.Range("K27").Value = 200
$ws.Range("M27").Value = -93
$ws.Range("H40").Value = 39970.71
$ws.Range("I40").Value = 48680.49
$ws.Range("K40").Value = 48680.49
$ws.Range("M40").Value = -48544.49
$ws.Range("H100").Value = 4423.364
$ws.Range("I100").Value = 3855.4285
$ws.Range("J100").Value = 5417.25
$ws.Range("K100").Value = 3855.4285
$ws.Range("L100").Value = 5417.25
$ws.Range("M100").Value = -3314.4285
$ws.Range("N100").Value = -6499.25
$ws.Range("H126").Value = 47298.15
$ws.Range("J126").Value = 2252
$ws.Range("L126").Value = 6756
$ws.Range("N126").Value = -11696
$ws.Range("H132").Value = 3060.4285
$ws.Range("I132").Value = 2895.4255
$ws.Range("K132").Value = 8686.2765
$ws.Range("M132").Value = -6156.2765

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 7961.5654
$ws.Range("J132").Value = 7896.4
$ws.Range("L132").Value = 23689.2
$ws.Range("N132").Value = -28749.2
$ws.Range("H136").Value = 2311.39
$ws.Range("I136").Value = 2364.7
$ws.Range("J136").Value = 2015.2222
$ws.Range("K136").Value = 7094.099999999999
$ws.Range("L136").Value = 6045.6666
$ws.Range("M136").Value = -4544.099999999999
$ws.Range("N136").Value = -11145.6666
